$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 12317.6

$ws.Range("H36").Value = 12317.6

$ws.Range("H54").Value = 3553.7144
$ws.Range("I54").Value = 2015.2
$ws.Range("K54").Value = 2015.2
$ws.Range("M54").Value = -1529.2

$ws.Range("H64").Value = 4140.607
$ws.Range("I64").Value = 3217
$ws.Range("J64").Value = 5568
$ws.Range("K64").Value = 3217
$ws.Range("L64").Value = 5568
$ws.Range("M64").Value = -2969
$ws.Range("N64").Value = -6064

$ws.Range("H67").Value = 4140.607
$ws.Range("I67").Value = 3217
$ws.Range("J67").Value = 5568
$ws.Range("K67").Value = 3217
$ws.Range("L67").Value = 5568
$ws.Range("M67").Value = -2359
$ws.Range("N67").Value = -7284

$ws.Range("H70").Value = 1907.9231
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 2036.6364
$ws.Range("K70").Value = 3600
$ws.Range("L70").Value = 6109.9092
$ws.Range("M70").Value = -3330
$ws.Range("N70").Value = -6649.9092

$ws.Range("H73").Value = 1907.9231
$ws.Range("I73").Value = 1200
$ws.Range("J73").Value = 2036.6364
$ws.Range("K73").Value = 3600
$ws.Range("L73").Value = 6109.9092
$ws.Range("M73").Value = -2664
$ws.Range("N73").Value = -7981.9092

$ws.Range("H92").Value = 427889.53
$ws.Range("I92").Value = 463505.34
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 463505.34
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = -462257.34
$ws.Range("N92").Value = -2996

$ws.Range("H116").Value = 9886292
$ws.Range("I116").Value = 13839639
$ws.Range("J116").Value = 2925
$ws.Range("K116").Value = 13839639
$ws.Range("L116").Value = 2925
$ws.Range("M116").Value = -13836197
$ws.Range("N116").Value = -9809

$ws.Range("H132").Value = 454050.94
$ws.Range("I132").Value = 578677.5600000001
$ws.Range("J132").Value = 17857.666
$ws.Range("K132").Value = 1736032.68
$ws.Range("L132").Value = 53572.99800000001
$ws.Range("M132").Value = -1733502.68
$ws.Range("N132").Value = -58632.99800000001

$ws.Range("H138").Value = 5955080
$ws.Range("I138").Value = 2509128
$ws.Range("J138").Value = 7410037.5
$ws.Range("K138").Value = 7527384
$ws.Range("L138").Value = 22230112.5
$ws.Range("M138").Value = -7522244
$ws.Range("N138").Value = -22240392.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 26334.334
$ws.Range("I6").Value = 38001.5
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 38001.5
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -37828.5
$ws.Range("N6").Value = -3346

$ws.Range("H32").Value = 26861.09
$ws.Range("I32").Value = 4184.1665
$ws.Range("J32").Value = 344338
$ws.Range("K32").Value = 4184.1665
$ws.Range("L32").Value = 344338
$ws.Range("M32").Value = -3897.1665
$ws.Range("N32").Value = -344912

$ws.Range("H63").Value = 35200
$ws.Range("I63").Value = 35200
$ws.Range("K63").Value = 35200
$ws.Range("M63").Value = -34514

$ws.Range("H66").Value = 35200
$ws.Range("I66").Value = 35200
$ws.Range("K66").Value = 176000
$ws.Range("M66").Value = -172568

$ws.Range("H74").Value = 6306.4585
$ws.Range("I74").Value = 917.8
$ws.Range("J74").Value = 33249.75
$ws.Range("K74").Value = 917.8
$ws.Range("L74").Value = 33249.75
$ws.Range("M74").Value = -43.79999999999995
$ws.Range("N74").Value = -34997.75

$ws.Range("H77").Value = 6306.4585
$ws.Range("I77").Value = 917.8
$ws.Range("J77").Value = 33249.75
$ws.Range("K77").Value = 4589
$ws.Range("L77").Value = 166248.75
$ws.Range("M77").Value = -221
$ws.Range("N77").Value = -174984.75

$ws.Range("H122").Value = 4361.25
$ws.Range("I122").Value = 1903.8889
$ws.Range("K122").Value = 5711.6667
$ws.Range("M122").Value = -3261.6667

$ws.Range("H132").Value = 2896.6667
$ws.Range("I132").Value = 2718.7073
$ws.Range("J132").Value = 3626.3
$ws.Range("K132").Value = 8156.1219
$ws.Range("L132").Value = 10878.9
$ws.Range("M132").Value = -5626.1219
$ws.Range("N132").Value = -15938.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3294.0417
$ws.Range("I134").Value = 2727.4736
$ws.Range("J134").Value = 5447
$ws.Range("K134").Value = 8182.4208
$ws.Range("L134").Value = 16341
$ws.Range("M134").Value = -5647.4208
$ws.Range("N134").Value = -21411

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 500250
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 1000000
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 1000000
$ws.Range("M13").Value = -361
$ws.Range("N13").Value = -1000278

$ws.Range("H31").Value = 3444.434
$ws.Range("I31").Value = 1224.44
$ws.Range("J31").Value = 5426.5713
$ws.Range("K31").Value = 1224.44
$ws.Range("L31").Value = 5426.5713
$ws.Range("M31").Value = -929.4400000000001
$ws.Range("N31").Value = -6016.5713

$ws.Range("H34").Value = 3444.434
$ws.Range("I34").Value = 1224.44
$ws.Range("J34").Value = 5426.5713
$ws.Range("K34").Value = 1224.44
$ws.Range("L34").Value = 5426.5713
$ws.Range("M34").Value = -1022.44
$ws.Range("N34").Value = -5830.5713

$ws.Range("H58").Value = 2838.037
$ws.Range("I58").Value = 1140.1818
$ws.Range("J58").Value = 4005.3125
$ws.Range("K58").Value = 1140.1818
$ws.Range("L58").Value = 4005.3125
$ws.Range("M58").Value = -937.1818000000001
$ws.Range("N58").Value = -4411.3125

$ws.Range("H62").Value = 16256.9375
$ws.Range("I62").Value = 26177.777
$ws.Range("K62").Value = 26177.777
$ws.Range("M62").Value = -25553.777

$ws.Range("H65").Value = 16256.9375
$ws.Range("I65").Value = 26177.777
$ws.Range("K65").Value = 130888.885
$ws.Range("M65").Value = -127768.885

$ws.Range("H134").Value = 3126.8125
$ws.Range("I134").Value = 1390.1111
$ws.Range("J134").Value = 5359.7144
$ws.Range("K134").Value = 4170.3333
$ws.Range("L134").Value = 16079.1432
$ws.Range("M134").Value = -1635.3333
$ws.Range("N134").Value = -21149.1432

$ws.Range("H136").Value = 2838.037
$ws.Range("I136").Value = 1140.1818
$ws.Range("J136").Value = 4005.3125
$ws.Range("K136").Value = 3420.5454
$ws.Range("L136").Value = 12015.9375
$ws.Range("M136").Value = -870.5454
$ws.Range("N136").Value = -17115.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 3687.5
$ws.Range("J105").Value = 3687.5
$ws.Range("L105").Value = 11062.5
$ws.Range("N105").Value = -16304.5

$ws.Range("H129").Value = 526
$ws.Range("I129").Value = 337.14285
$ws.Range("J129").Value = 966.6667
$ws.Range("K129").Value = 1011.42855
$ws.Range("L129").Value = 2900.0001
$ws.Range("M129").Value = 3988.57145
$ws.Range("N129").Value = -12900.0001

$ws.Range("H131").Value = 3832593.2
$ws.Range("J131").Value = 4066275
$ws.Range("L131").Value = 12198825
$ws.Range("N131").Value = -12208905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2286.8438
$ws.Range("I102").Value = 2167.3462
$ws.Range("J102").Value = 2804.6667
$ws.Range("K102").Value = 2167.3462
$ws.Range("L102").Value = 2804.6667
$ws.Range("M102").Value = -545.3462
$ws.Range("N102").Value = -6048.6667

$ws.Range("H122").Value = 655603.4
$ws.Range("I122").Value = 1112311.1
$ws.Range("K122").Value = 3336933.3
$ws.Range("M122").Value = -3334483.3

$ws.Range("H123").Value = 10953.053
$ws.Range("J123").Value = 10953.053
$ws.Range("L123").Value = 10953.053
$ws.Range("N123").Value = -15853.053

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2056.7144
$ws.Range("I68").Value = 1849.8334
$ws.Range("K68").Value = 1849.8334
$ws.Range("M68").Value = -1100.8334

$ws.Range("H71").Value = 2056.7144
$ws.Range("I71").Value = 1849.8334
$ws.Range("K71").Value = 9249.166999999999
$ws.Range("M71").Value = -5505.166999999999

$ws.Range("H122").Value = 3585.3076
$ws.Range("I122").Value = 2876
$ws.Range("J122").Value = 3900.5557
$ws.Range("K122").Value = 8628
$ws.Range("L122").Value = 11701.6671
$ws.Range("M122").Value = -6178
$ws.Range("N122").Value = -16601.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 57170000
$ws.Range("J2").Value = 70003
$ws.Range("L2").Value = 70003
$ws.Range("N2").Value = -70227

$ws.Range("H11").Value = 50000
$ws.Range("I11").Value = 50000
$ws.Range("K11").Value = 50000
$ws.Range("M11").Value = -49858

$ws.Range("H100").Value = 257
$ws.Range("I100").Value = 158.83333
$ws.Range("J100").Value = 551.5
$ws.Range("K100").Value = 317.66666
$ws.Range("L100").Value = 1103
$ws.Range("M100").Value = 223.33334
$ws.Range("N100").Value = -2185

$ws.Range("H122").Value = 92068.17999999999
$ws.Range("I122").Value = 100955
$ws.Range("K122").Value = 302865
$ws.Range("M122").Value = -300415

$ws.Range("H132").Value = 2227.9333
$ws.Range("I132").Value = 2007.0555
$ws.Range("J132").Value = 3111.4443
$ws.Range("K132").Value = 6021.166499999999
$ws.Range("L132").Value = 9334.332900000001
$ws.Range("M132").Value = -3491.166499999999
$ws.Range("N132").Value = -14394.3329

$ws.Range("H136").Value = 1717.0731
$ws.Range("I136").Value = 1185.7826
$ws.Range("J136").Value = 2395.9443
$ws.Range("K136").Value = 4005.3478
$ws.Range("L136").Value = 7187.8329
$ws.Range("M136").Value = -1007.3478
$ws.Range("N136").Value = -12287.8329
